# "expanded neighborhood to include not feasible"
# Update the benchmark "instance0X" results (column I, row 4-10) with the
# refreshed run numbers, underline the run-date cell (I2), and leave the
# selection on the first re-measured cell (I5) the way the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Refreshed benchmark values (column I) for instance01..instance07.
$ws.Range("I4").Value  = 169.42500000000001
$ws.Range("I5").Value  = 45.341000000000001
$ws.Range("I6").Value  = 49.561
$ws.Range("I7").Value  = 11.500999999999999
$ws.Range("I8").Value  = 18.212
$ws.Range("I9").Value  = 4.6239999999999997
$ws.Range("I10").Value = 12.855

# Mark the run date (I2) with an underline, matching the author's edit.
$ws.Range("I2").Font.Underline = $true

# Leave the selection where the author left it after the edit.
$ws.Range("I5").Select()

$wb.Application.Calculate()
